$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set header cells I1 and J1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header formatting already used by the rest of row 1 (e.g. H1: bold,
# bordered, centered) by copying H1's format onto the two new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Set the new data cells I2 and J2
$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 7
